$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19 (student #16, Овечкин Родион): columns G:J (ДЗ_5, ДЗ_6, ДЗ_7, Лаб_1) 0 -> 5 ---
$ws.Range("G19").Value = 5
$ws.Range("H19").Value = 5
$ws.Range("I19").Value = 5
$ws.Range("J19").Value = 5
$ws.Range("O19").Value = "BPV"

# --- Row 23 (student #20, Уваров Арсений): columns G:J (ДЗ_5, ДЗ_6, ДЗ_7, Лаб_1) 0 -> 5 ---
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = 5
$ws.Range("I23").Value = 5
$ws.Range("J23").Value = 5
$ws.Range("O23").Value = "BPV"

# --- Update the active selection to match where the author ended up ---
$ws.Range("O24").Select()

$wb.Save()
